$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y (25) and AA (27) hold date-looking text ("Startdatum"/"Slutdatum").
# Force the cell to Text format before assignment, otherwise Excel would
# auto-convert the string into a date serial number.

# ---- Row 26 ----
$ws.Cells.Item(26, 1).Value = 112093192
$ws.Cells.Item(26, 2).Value = 90687
$ws.Cells.Item(26, 3).Value = "Ovaliderad"
$ws.Cells.Item(26, 4).Value = "LC"
$ws.Cells.Item(26, 5).Value = 5964
$ws.Cells.Item(26, 6).Value = "Fjällig taggsvamp s.str."
$ws.Cells.Item(26, 7).Value = "Sarcodon imbricatus s.str."
$ws.Cells.Item(26, 8).Value = "(L.:Fr.) P.Karst."
$ws.Cells.Item(26, 16).Value = "Fagerdal, Upl"
$ws.Cells.Item(26, 17).Value = 639179.9128251362
$ws.Cells.Item(26, 18).Value = 6701165.391882338
$ws.Cells.Item(26, 19).Value = 10
$ws.Cells.Item(26, 20).Value = "Uppsala"
$ws.Cells.Item(26, 21).Value = "Tierp"
$ws.Cells.Item(26, 22).Value = "Uppland"
$ws.Cells.Item(26, 23).Value = "Tolfta"
$ws.Cells.Item(26, 25).NumberFormat = "@"
$ws.Cells.Item(26, 25).Value = "2023-09-12"
$ws.Cells.Item(26, 26).Value = "00:00"
$ws.Cells.Item(26, 27).NumberFormat = "@"
$ws.Cells.Item(26, 27).Value = "2023-09-12"
$ws.Cells.Item(26, 28).Value = "00:00"
$ws.Cells.Item(26, 30).Value = $false
$ws.Cells.Item(26, 31).Value = $false
$ws.Cells.Item(26, 33).Value = $false
$ws.Cells.Item(26, 49).Value = "Samuel Persson"
$ws.Cells.Item(26, 50).Value = "Samuel Persson"

# ---- Row 27 ----
$ws.Cells.Item(27, 1).Value = 112093193
$ws.Cells.Item(27, 2).Value = 89183
$ws.Cells.Item(27, 3).Value = "Ovaliderad"
$ws.Cells.Item(27, 4).Value = "LC"
$ws.Cells.Item(27, 5).Value = 3215
$ws.Cells.Item(27, 6).Value = "Rödgul trumpetsvamp"
$ws.Cells.Item(27, 7).Value = "Craterellus lutescens"
$ws.Cells.Item(27, 8).Value = "(Fr.) Fr."
$ws.Cells.Item(27, 16).Value = "Fagerdal, Upl"
$ws.Cells.Item(27, 17).Value = 639179.9128251362
$ws.Cells.Item(27, 18).Value = 6701165.391882338
$ws.Cells.Item(27, 19).Value = 10
$ws.Cells.Item(27, 20).Value = "Uppsala"
$ws.Cells.Item(27, 21).Value = "Tierp"
$ws.Cells.Item(27, 22).Value = "Uppland"
$ws.Cells.Item(27, 23).Value = "Tolfta"
$ws.Cells.Item(27, 25).NumberFormat = "@"
$ws.Cells.Item(27, 25).Value = "2023-09-12"
$ws.Cells.Item(27, 26).Value = "00:00"
$ws.Cells.Item(27, 27).NumberFormat = "@"
$ws.Cells.Item(27, 27).Value = "2023-09-12"
$ws.Cells.Item(27, 28).Value = "00:00"
$ws.Cells.Item(27, 30).Value = $false
$ws.Cells.Item(27, 31).Value = $false
$ws.Cells.Item(27, 33).Value = $false
$ws.Cells.Item(27, 49).Value = "Samuel Persson"
$ws.Cells.Item(27, 50).Value = "Samuel Persson"

# ---- Row 28 ----
$ws.Cells.Item(28, 1).Value = 112093190
$ws.Cells.Item(28, 2).Value = 85210
$ws.Cells.Item(28, 3).Value = "Ovaliderad"
$ws.Cells.Item(28, 4).Value = "LC"
$ws.Cells.Item(28, 5).Value = 3624
$ws.Cells.Item(28, 6).Value = "Strimspindling"
$ws.Cells.Item(28, 7).Value = "Cortinarius glaucopus"
$ws.Cells.Item(28, 8).Value = "(Schaeff. : Fr.) Fr."
$ws.Cells.Item(28, 16).Value = "Fagerdal, Upl"
$ws.Cells.Item(28, 17).Value = 639179.9128251362
$ws.Cells.Item(28, 18).Value = 6701165.391882338
$ws.Cells.Item(28, 19).Value = 10
$ws.Cells.Item(28, 20).Value = "Uppsala"
$ws.Cells.Item(28, 21).Value = "Tierp"
$ws.Cells.Item(28, 22).Value = "Uppland"
$ws.Cells.Item(28, 23).Value = "Tolfta"
$ws.Cells.Item(28, 25).NumberFormat = "@"
$ws.Cells.Item(28, 25).Value = "2023-09-12"
$ws.Cells.Item(28, 26).Value = "00:00"
$ws.Cells.Item(28, 27).NumberFormat = "@"
$ws.Cells.Item(28, 27).Value = "2023-09-12"
$ws.Cells.Item(28, 28).Value = "00:00"
$ws.Cells.Item(28, 30).Value = $false
$ws.Cells.Item(28, 31).Value = $false
$ws.Cells.Item(28, 33).Value = $false
$ws.Cells.Item(28, 49).Value = "Samuel Persson"
$ws.Cells.Item(28, 50).Value = "Samuel Persson"

# ---- Row 29 ----
$ws.Cells.Item(29, 1).Value = 112093186
$ws.Cells.Item(29, 2).Value = 88909
$ws.Cells.Item(29, 3).Value = "Ovaliderad"
$ws.Cells.Item(29, 4).Value = "VU"
$ws.Cells.Item(29, 5).Value = 720
$ws.Cells.Item(29, 6).Value = "Violgubbe"
$ws.Cells.Item(29, 7).Value = "Gomphus clavatus"
$ws.Cells.Item(29, 8).Value = "(Pers.) Gray"
$ws.Cells.Item(29, 16).Value = "Fagerdal, Upl"
$ws.Cells.Item(29, 17).Value = 639204.9761395331
$ws.Cells.Item(29, 18).Value = 6701015.582563667
$ws.Cells.Item(29, 19).Value = 10
$ws.Cells.Item(29, 20).Value = "Uppsala"
$ws.Cells.Item(29, 21).Value = "Tierp"
$ws.Cells.Item(29, 22).Value = "Uppland"
$ws.Cells.Item(29, 23).Value = "Tolfta"
$ws.Cells.Item(29, 25).NumberFormat = "@"
$ws.Cells.Item(29, 25).Value = "2023-09-12"
$ws.Cells.Item(29, 26).Value = "00:00"
$ws.Cells.Item(29, 27).NumberFormat = "@"
$ws.Cells.Item(29, 27).Value = "2023-09-12"
$ws.Cells.Item(29, 28).Value = "00:00"
$ws.Cells.Item(29, 29).Value = "Till stor del barkborredödat bestånd intill stort kalhygge"
$ws.Cells.Item(29, 30).Value = $true
$ws.Cells.Item(29, 31).Value = $false
$ws.Cells.Item(29, 33).Value = $false
$ws.Cells.Item(29, 49).Value = "Samuel Persson"
$ws.Cells.Item(29, 50).Value = "Samuel Persson"
